# "fixed errors in Monte Carlo with addition of discrete measurements,
#  documented results"
#
# - general!B8  (n_MonteCarloRuns)      20 -> 200
# - general!B13 (process_GPS_enable)     0 -> 1   (enable the discrete GPS
#   measurement update)
# - remove the now-unused processNoiseOn / measurementNoiseOn rows (18:19)
#   from the "general" sheet entirely
# - leave the cursor on general!C13 and truthStateParams!D42, matching the
#   author's final selection when the workbook was saved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Activate()

# number of Monte Carlo runs: 20 -> 200
$ws.Range("B8").Value = 200

# process_GPS_enable: 0 -> 1
$ws.Range("B13").Value = 1

# delete rows 18 ("measurementNoiseOn") and 19 ("processNoiseOn") entirely -
# everything below shifts up two rows
$ws.Rows("18:19").Delete()

# author's final selection on the general sheet
$ws.Range("C13").Select()

# author's final selection on truthStateParams
$ws2 = $wb.Worksheets.Item("truthStateParams")
$ws2.Activate()
$ws2.Range("D42").Select()

# general stays the tab that is shown/selected when the file is reopened
$ws.Activate()
